$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'61.493.41"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -3.69%  "
$ws.Range("D3").Value = "'2.465.94"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -6.41%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").Value = "'551.33"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -4.69%  "
$ws.Range("D6").Value = "'147.03"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -6.32%  "
$ws.Range("D7").Value = "'0.999"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.09%  "
$ws.Range("D8").Value = "'0.593"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -6.01%  "
$ws.Range("D9").Value = "'2.467.76"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -6.27%  "
$ws.Range("E10").Value = "  -9.20%  "
$ws.Range("D11").Value = "'5.46"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -6.49%  "
$ws.Range("E12").Value = "  -1.69%  "
$ws.Range("D13").Value = "'0.354"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -7.76%  "
$ws.Range("D14").Value = "'26.15"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -8.78%  "
$ws.Range("D15").Value = "'2.905.81"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -6.56%  "
$ws.Range("D16").Value = "'0.0000167"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -9.20%  "
$ws.Range("D17").Value = "'61.394.30"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -3.70%  "
$ws.Range("D18").Value = "'2.462.24"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -6.53%  "
$ws.Range("D19").Value = "'11.14"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -8.41%  "
$ws.Range("D20").Value = "'7.14"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -8.03%  "
$ws.Range("D21").Value = "'4.20"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -7.15%  "
$ws.Range("D22").Value = "'318.75"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -7.28%  "
$ws.Range("D24").Value = "'1.89"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.58%  "
$ws.Range("D25").Value = "'64.00"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -6.42%  "
$ws.Range("D26").Value = "'0.0₃0986"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -12.55%  "
$ws.Range("B27").Value = "Bittensor"
$ws.Range("C27").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D27").Value = "'553.54"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -4.81%  "
$ws.Range("B28").Value = "WrappedeETH"
$ws.Range("C28").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("D28").Value = "'2.580.68"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -6.76%  "
$ws.Range("E29").Value = "  +0.04%  "
$ws.Range("D30").Value = "'1.48"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -9.38%  "
$ws.Range("D31").Value = "'8.32"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -10.12%  "
$ws.Range("D32").Value = "'7.70"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -6.42%  "
$ws.Range("E33").Value = "  -8.63%  "
$ws.Range("D34").Value = "'1.91"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -7.33%  "
$ws.Range("D35").Value = "'1.60"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -7.80%  "
$ws.Range("D36").Value = "'5.92"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -10.61%  "
$ws.Range("B37").Value = "NEARProtocol"
$ws.Range("C37").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D37").Value = "'4.86"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -10.58%  "
$ws.Range("B38").Value = "FirstDigitalUSD"
$ws.Range("C38").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D38").Value = "'0.998"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.09%  "
$ws.Range("D39").Value = "'0.381"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -5.57%  "
$ws.Range("D40").Value = "'18.43"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -6.91%  "
$ws.Range("B41").Value = "Stacks"
$ws.Range("C41").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D41").Value = "'1.78"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -6.90%  "
$ws.Range("B42").Value = "Monero"
$ws.Range("C42").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D42").Value = "'142.39"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -7.08%  "
$ws.Range("D44").Value = "'40.44"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -4.12%  "
$ws.Range("E45").Value = "  -7.28%  "
$ws.Range("D46").Value = "'146.64"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -9.63%  "
$ws.Range("D47").Value = "'3.62"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -7.78%  "
$ws.Range("D48").Value = "'21.63"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -10.45%  "
$ws.Range("D49").Value = "'0.0539"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -8.41%  "
$ws.Range("D50").Value = "'0.592"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -6.47%  "
$ws.Range("D51").Value = "'0.0939"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -6.47%  "
